# --- edit.ps1 ---
# Restructure the "물건" (items) sheet into the new ObjID/ObjName/ObjInfo/defaultPrice/expensive/tooExpensive schema
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Drop the old stray accessory block (potion/8/.../2) that lived in H6:K6 -- it is being
# promoted into a full row (row 9) of the main table instead.
$ws.Range("H6:K6").ClearContents()

# Row 1
$ws.Range("A1").Value = "ObjID"
$ws.Range("B1").Value = "ObjName"
$ws.Range("C1").Value = "ObjInfo"
$ws.Range("D1").Value = "defaultPrice"
$ws.Range("E1").Value = "expensive"
$ws.Range("F1").Value = "tooExpensive"

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "pumpkin"
$ws.Range("C2").Value = "큼지막하고 맛있어보이는 호박이다. 다른 호박이랑 다른 점을 모르겠다."
$ws.Range("D2").Value = 24
$ws.Range("E2").Value = 30
$ws.Range("F2").Value = 50

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "cake"
$ws.Range("C3").Value = "냉동 딸기 쇼트 케이크. 다른 케이크와 다른 점은 없어 보인다. 케이크가 신선해보이는 건 거짓말이다. 누가봐도 공장에서 나온 케이크다."
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 20

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "door"
$ws.Range("C4").Value = "유리창이 있는 하얀색 문"
$ws.Range("D4").Value = 65
$ws.Range("E4").Value = 90
$ws.Range("F4").Value = 150

# Row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "hoodie"
$ws.Range("C5").Value = "빨간색 후드티. 팔쪽에 흰색과 검은색으로 포인트 디자인이 있다."
$ws.Range("D5").Value = 43
$ws.Range("E5").Value = 80
$ws.Range("F5").Value = 100

# Row 6
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "cat_tower"
$ws.Range("C6").Value = "3단으로 된 캣타워. 핑크색 포인트 색깔이 인상적이다."
$ws.Range("D6").Value = 72
$ws.Range("E6").Value = 100
$ws.Range("F6").Value = 150

# Row 7
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "piano"
$ws.Range("C7").Value = "작은 피아노. 원래 피아노가 비싼 편이라 사려면 큰 마음을 먹어야한다."
$ws.Range("D7").Value = 700
$ws.Range("E7").Value = 900
$ws.Range("F7").Value = 1200

# Row 8
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "camera"
$ws.Range("C8").Value = "디지털 카메라다. 얼마나 비싼지 모르겠다."
$ws.Range("D8").Value = 560
$ws.Range("E8").Value = 780
$ws.Range("F8").Value = 1000

# Row 9
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "potion"
$ws.Range("C9").Value = "수상하게 생긴 포션. 마치 판타지 세상에서 튀어나온 것 같이 생겼다. 얼마인지 모르겠다. 효능에 따라 아무 효능이 없으면 2크레딧, 효능이 있으면 50크레딧정도 할것 같다."
$ws.Range("D9").Value = 20
$ws.Range("E9").Value = 35
$ws.Range("F9").Value = 60

